$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 2650.195
$ws.Range("I80").Value = 1197.9166
$ws.Range("J80").Value = 3251.138
$ws.Range("K80").Value = 3593.7498
$ws.Range("L80").Value = 9753.414000000001
$ws.Range("M80").Value = -2595.7498
$ws.Range("N80").Value = -11749.414
# Row 83
$ws.Range("H83").Value = 2650.195
$ws.Range("I83").Value = 1197.9166
$ws.Range("J83").Value = 3251.138
$ws.Range("K83").Value = 10781.2494
$ws.Range("L83").Value = 29260.242
$ws.Range("M83").Value = -5789.249400000001
$ws.Range("N83").Value = -39244.242
# Row 137
$ws.Range("H137").Value = 2285.6553
$ws.Range("I137").Value = 1702.9524
$ws.Range("J137").Value = 3815.25
$ws.Range("K137").Value = 5108.857199999999
$ws.Range("L137").Value = 11445.75
$ws.Range("M137").Value = -2558.857199999999
$ws.Range("N137").Value = -16545.75
# Row 138
$ws.Range("H138").Value = 3208.6099
$ws.Range("J138").Value = 3865.3794
$ws.Range("L138").Value = 11596.1382
$ws.Range("N138").Value = -21876.1382

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13990.47
$ws.Range("I32").Value = 10514.779
$ws.Range("K32").Value = 10514.779
$ws.Range("M32").Value = -10227.779
# Row 61
$ws.Range("H61").Value = 5806.8335
$ws.Range("I61").Value = 4502.3335
$ws.Range("K61").Value = 4502.3335
$ws.Range("M61").Value = -4290.3335
# Row 88
$ws.Range("H88").Value = 4637291
$ws.Range("I88").Value = 19380.666
$ws.Range("J88").Value = 6946246
$ws.Range("K88").Value = 19380.666
$ws.Range("L88").Value = 6946246
$ws.Range("M88").Value = -18974.666
$ws.Range("N88").Value = -6947058
# Row 91
$ws.Range("H91").Value = 4637291
$ws.Range("I91").Value = 19380.666
$ws.Range("J91").Value = 6946246
$ws.Range("K91").Value = 19380.666
$ws.Range("L91").Value = 6946246
$ws.Range("M91").Value = -17976.666
$ws.Range("N91").Value = -6949054
# Row 125
$ws.Range("H125").Value = 70178.75
$ws.Range("J125").Value = 70178.75
$ws.Range("L125").Value = 70178.75
$ws.Range("N125").Value = -80018.75
# Row 136
$ws.Range("H136").Value = 5806.8335
$ws.Range("I136").Value = 4502.3335
$ws.Range("K136").Value = 13507.0005
$ws.Range("M136").Value = -10957.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4166
$ws.Range("I20").Value = 4001.2
$ws.Range("K20").Value = 4001.2
$ws.Range("M20").Value = -3754.2
# Row 86
$ws.Range("H86").Value = 7588.6665
$ws.Range("I86").Value = 5048.5
$ws.Range("J86").Value = 12669
$ws.Range("K86").Value = 5048.5
$ws.Range("L86").Value = 12669
$ws.Range("M86").Value = -3925.5
$ws.Range("N86").Value = -14915
# Row 89
$ws.Range("H89").Value = 7588.6665
$ws.Range("I89").Value = 5048.5
$ws.Range("J89").Value = 12669
$ws.Range("K89").Value = 25242.5
$ws.Range("L89").Value = 63345
$ws.Range("M89").Value = -19626.5
$ws.Range("N89").Value = -74577
# Row 94
$ws.Range("H94").Value = 10002386
$ws.Range("I94").Value = 4763912
$ws.Range("J94").Value = 22225494
$ws.Range("K94").Value = 4763912
$ws.Range("L94").Value = 22225494
$ws.Range("M94").Value = -4763461
$ws.Range("N94").Value = -22226396
# Row 99
$ws.Range("H99").Value = 23013.445
$ws.Range("I99").Value = 24131.234
$ws.Range("J99").Value = 4011
$ws.Range("K99").Value = 24131.234
$ws.Range("L99").Value = 4011
$ws.Range("M99").Value = -22633.234
$ws.Range("N99").Value = -7007
# Row 105
$ws.Range("H105").Value = 4044.8635
$ws.Range("I105").Value = 2157.5
$ws.Range("K105").Value = 2157.5
$ws.Range("M105").Value = -410.5
# Row 107
$ws.Range("H107").Value = 2624.5
$ws.Range("I107").Value = 2624.5
$ws.Range("K107").Value = 2624.5
$ws.Range("M107").Value = -704.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3453.0364
$ws.Range("I31").Value = 2339.7273
$ws.Range("J31").Value = 4195.242
$ws.Range("K31").Value = 2339.7273
$ws.Range("L31").Value = 4195.242
$ws.Range("M31").Value = -2044.7273
$ws.Range("N31").Value = -4785.242
# Row 34
$ws.Range("H34").Value = 3453.0364
$ws.Range("I34").Value = 2339.7273
$ws.Range("J34").Value = 4195.242
$ws.Range("K34").Value = 2339.7273
$ws.Range("L34").Value = 4195.242
$ws.Range("M34").Value = -2137.7273
$ws.Range("N34").Value = -4599.242
# Row 99
$ws.Range("H99").Value = 24939716
$ws.Range("I99").Value = 6106486
$ws.Range("J99").Value = 40006300
$ws.Range("K99").Value = 6106486
$ws.Range("L99").Value = 40006300
$ws.Range("M99").Value = -6104988
$ws.Range("N99").Value = -40009296
# Row 103
$ws.Range("H103").Value = 12481.857
$ws.Range("I103").Value = 12481.857
$ws.Range("K103").Value = 12481.857
$ws.Range("M103").Value = -11309.857
# Row 122
$ws.Range("H122").Value = 4657297.5
$ws.Range("I122").Value = 9310183
$ws.Range("J122").Value = 4411.3687
$ws.Range("K122").Value = 27930549
$ws.Range("L122").Value = 13234.1061
$ws.Range("M122").Value = -27928099
$ws.Range("N122").Value = -18134.1061
# Row 126
$ws.Range("H126").Value = 24939716
$ws.Range("I126").Value = 6106486
$ws.Range("J126").Value = 40006300
$ws.Range("K126").Value = 18319458
$ws.Range("L126").Value = 120018900
$ws.Range("M126").Value = -18316988
$ws.Range("N126").Value = -120023840
# Row 134
$ws.Range("H134").Value = 3215.8076
$ws.Range("I134").Value = 1780.55
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 5341.65
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -2806.65
$ws.Range("N134").Value = -29070

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 15626158
$ws.Range("I97").Value = 919.1
$ws.Range("J97").Value = 41668224
$ws.Range("K97").Value = 919.1
$ws.Range("L97").Value = 41668224
$ws.Range("M97").Value = -423.1
$ws.Range("N97").Value = -41669216
# Row 107
$ws.Range("H107").Value = 276.92307
$ws.Range("I107").Value = 270
$ws.Range("J107").Value = 315
$ws.Range("K107").Value = 270
$ws.Range("L107").Value = 315
$ws.Range("M107").Value = 1650
$ws.Range("N107").Value = -4155
# Row 132
$ws.Range("H132").Value = 4863.2666
$ws.Range("I132").Value = 3992.125
$ws.Range("K132").Value = 11976.375
$ws.Range("M132").Value = -9446.375
# Row 135
$ws.Range("H135").Value = 61249.875
$ws.Range("J135").Value = 61249.875
$ws.Range("L135").Value = 61249.875
$ws.Range("N135").Value = -71389.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 7133.0557
$ws.Range("J68").Value = 7441.6665
$ws.Range("L68").Value = 7441.6665
$ws.Range("N68").Value = -8939.666499999999
# Row 71
$ws.Range("H71").Value = 7133.0557
$ws.Range("J71").Value = 7441.6665
$ws.Range("L71").Value = 37208.3325
$ws.Range("N71").Value = -44696.3325
# Row 100
$ws.Range("H100").Value = 116178.1
$ws.Range("I100").Value = 371795.34
$ws.Range("J100").Value = 6627.857
$ws.Range("K100").Value = 371795.34
$ws.Range("L100").Value = 6627.857
$ws.Range("M100").Value = -371254.34
$ws.Range("N100").Value = -7709.857

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 11758.667
$ws.Range("I74").Value = 9824
$ws.Range("J74").Value = 12726
$ws.Range("K74").Value = 9824
$ws.Range("L74").Value = 12726
$ws.Range("M74").Value = -8888
$ws.Range("N74").Value = -14598
# Row 75
$ws.Range("H75").Value = 67976.336
$ws.Range("J75").Value = 67976.336
$ws.Range("L75").Value = 67976.336
$ws.Range("N75").Value = -69848.336
# Row 77
$ws.Range("H77").Value = 11758.667
$ws.Range("I77").Value = 9824
$ws.Range("J77").Value = 12726
$ws.Range("K77").Value = 29472
$ws.Range("L77").Value = 38178
$ws.Range("M77").Value = -24792
$ws.Range("N77").Value = -47538
# Row 78
$ws.Range("H78").Value = 67976.336
$ws.Range("J78").Value = 67976.336
$ws.Range("L78").Value = 203929.008
$ws.Range("N78").Value = -213289.008
# Row 136
$ws.Range("H136").Value = 4367.794
$ws.Range("I136").Value = 3808.4
$ws.Range("K136").Value = 11425.2
$ws.Range("M136").Value = -8875.200000000001
